# The "Förändrad" (Changed) column (column C) is updated by one day for
# every data row: 2023-10-03 (serial 45202) -> 2023-10-04 (serial 45203).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (robust against API quirks).
$lastRowUsed = $ws.UsedRange.Rows.Count
$lastRowXlUp = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

$lastRow = $lastRowUsed
if ($lastRowXlUp -gt $lastRow) {
    $lastRow = $lastRowXlUp
}
if ($lastRow -lt 2) {
    $lastRow = 2
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
